$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'315.03"
$ws.Range("D3").Value = "'39.28"
$ws.Range("E3").Value = "'-0.49%"
$ws.Range("D4").Value = "'5.154"
$ws.Range("E4").Value = "'0.39%"
$ws.Range("D5").Value = "'0.08158"
$ws.Range("E5").Value = "'0.33%"
$ws.Range("D6").Value = "'1.972"
$ws.Range("E6").Value = "'1.33%"
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'4.411"
$ws.Range("E7").Value = "'4.37%"
$ws.Range("B8").Value = "KuCoinToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D8").Value = "'8.330"
$ws.Range("E8").Value = "'2.47%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").Value = "'0.9372"
$ws.Range("E9").Value = "'1.15%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1316"
$ws.Range("E10").Value = "'-6.96%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1986"
$ws.Range("E11").Value = "'2.95%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.09028"
$ws.Range("E12").Value = "'-0.48%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03489"
$ws.Range("E13").Value = "'-0.14%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09711"
$ws.Range("E14").Value = "'-0.98%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001411"
$ws.Range("E15").Value = "'1.39%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006071"
$ws.Range("E16").Value = "'2.76%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.610"
$ws.Range("E17").Value = "'-7.76%"
$ws.Range("D18").Value = "'3.191"
$ws.Range("E18").Value = "'-5.49%"
$ws.Range("E19").Value = "'0.43%"
$ws.Range("E20").Value = "'-0.31%"
$ws.Range("D21").Value = "'5.014"
$ws.Range("E21").Value = "'5.98%"
$ws.Range("D23").Value = "'0.04366"
$ws.Range("E23").Value = "'-0.50%"
$ws.Range("D24").Value = "'0.001244"
$ws.Range("E24").Value = "'0.97%"
$ws.Range("D25").Value = "'0.004747"
$ws.Range("E25").Value = "'-2.39%"
$ws.Range("D26").Value = "'0.0003895"
$ws.Range("E26").Value = "'199.27%"
$ws.Range("E27").Value = "'-7.70%"
$ws.Range("D39").Value = "'0.02239"
$ws.Range("E39").Value = "'8.30%"
$ws.Range("D40").Value = "'0.05247"
$ws.Range("E40").Value = "'3.57%"
$ws.Range("D41").Value = "'0.007550"
$ws.Range("E41").Value = "'1.53%"
$ws.Range("E42").Value = "'5.43%"
$ws.Range("D43").Value = "'0.1396"
$ws.Range("E43").Value = "'2.18%"
$ws.Range("D44").Value = "'0.002103"
$ws.Range("E44").Value = "'-1.40%"
$ws.Range("D45").Value = "'0.009139"
$ws.Range("E45").Value = "'-4.43%"
$ws.Range("D46").Value = "'0.00006830"
$ws.Range("E46").Value = "'7.06%"
$ws.Range("E47").Value = "'0.07%"
$ws.Range("D48").Value = "'0.003015"
$ws.Range("E48").Value = "'10.94%"
$ws.Range("E50").Value = "'0.07%"
$ws.Range("D51").Value = "'0.0002004"
$ws.Range("E51").Value = "'0.07%"
